$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: clear Common_name (C) and Category (D)
$ws.Range("C23:D23").ClearContents()

# Row 29: capitalize "Unassigned" in Species_name (B), clear C/D
$ws.Range("B29").Value = "Unassigned"
$ws.Range("C29:D29").ClearContents()

# Row 41: capitalize "Unassigned" in Species_name (B), clear C/D
$ws.Range("B41").Value = "Unassigned"
$ws.Range("C41:D41").ClearContents()

# Row 42 & 43: swap ASV_ID/Species/Common/Category content between the two rows
$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "Unassigned"
$ws.Range("C42:D42").ClearContents()
$ws.Range("J42").Value = 0

$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"
$ws.Range("J43").ClearContents()

# Row 55 & 56: swap ASV_ID/Species/Common/Category content between the two rows
$ws.Range("A55").Value = "975b1dbdc7405f6e27bf63893e91e0ed"
$ws.Range("B55").Value = "Centropristis striata"
$ws.Range("C55").Value = "Black sea bass"
$ws.Range("D55").Value = "Teleost Fish"

$ws.Range("A56").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B56").Value = "Unassigned"
$ws.Range("C56:D56").ClearContents()

# Row 60: capitalize "Unassigned" in Species_name (B), clear C/D
$ws.Range("B60").Value = "Unassigned"
$ws.Range("C60:D60").ClearContents()
